$wb = $excel.ActiveWorkbook

# The status text "Ready for handoff" is a single shared string referenced by
# several cells (Overview!E3/F3 "zh-cn"/"de-de" status columns, and the
# "Status" column (C3) on both the zh-cn and de-de detail sheets). The report
# generator flips all of them to "Handback transform failed" for the
# "36912af3-f468-42ce-9bed-7be3d6499204" row at once.
$newStatus = "Handback transform failed"

# --- Overview sheet: update status text for the "36912af3..." row ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# --- zh-cn sheet: same status update, plus widen Error Detail column (P) and
#     set the handback mismatch error message in P3 ---
# Note: the engine stores width = ColumnWidth + 5/6, so to land on a stored
# width of exactly 40 we request 40 - 5/6.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664
$wsZhCn.Range("P3").Value = "Handback file name: t0o0mdez.2nw is different with handoff file name: 36912af3-f468-42ce-9bed-7be3d6499204.e19917b9f62cbec26f137fc0b4587af2c379a7ec.zh-cn."

# --- de-de sheet: same status update, plus widen Error Detail column (P) and
#     set the handback mismatch error message in P3 ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664
$wsDeDe.Range("P3").Value = "Handback file name: t0o0mdez.2nw is different with handoff file name: 36912af3-f468-42ce-9bed-7be3d6499204.e19917b9f62cbec26f137fc0b4587af2c379a7ec.de-de."
